{"js": "// Update the two-digit multiplication problems in the table to the new\n// set of values (see commit \"Update master to output generated at aa3dc9e\").\nconst replacements = [\n  [\"23\u00d738=\", \"54\u00d738=\"],\n  [\"51\u00d740=\", \"88\u00d740=\"],\n  [\"91\u00d711=\", \"55\u00d769=\"],\n  [\"18\u00d756=\", \"21\u00d745=\"],\n  [\"51\u00d742=\", \"15\u00d785=\"],\n  [\"32\u00d788=\", \"56\u00d744=\"],\n  [\"51\u00d784=\", \"38\u00d790=\"],\n  [\"95\u00d756=\", \"97\u00d772=\"],\n  [\"43\u00d757=\", \"74\u00d735=\"],\n  [\"77\u00d715=\", \"20\u00d765=\"],\n  [\"17\u00d720=\", \"89\u00d737=\"],\n  [\"91\u00d762=\", \"86\u00d791=\"],\n  [\"56\u00d723=\", \"53\u00d744=\"],\n  [\"38\u00d715=\", \"64\u00d713=\"],\n  [\"59\u00d749=\", \"64\u00d712=\"],\n  [\"46\u00d754=\", \"29\u00d788=\"],\n  [\"57\u00d785=\", \"87\u00d738=\"],\n  [\"80\u00d735=\", \"34\u00d720=\"],\n  [\"94\u00d781=\", \"35\u00d764=\"],\n  [\"77\u00d729=\", \"41\u00d713=\"],\n  [\"41\u00d719=\", \"57\u00d711=\"],\n  [\"43\u00d749=\", \"11\u00d798=\"],\n  [\"48\u00d712=\", \"66\u00d725=\"],\n  [\"40\u00d775=\", \"62\u00d742=\"],\n  [\"82\u00d744=\", \"59\u00d743=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the two-digit multiplication problems in the table to the new\n# set of values (see commit \"Update master to output generated at aa3dc9e\").\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"23\u00d738=\"; New = \"54\u00d738=\" },\n    @{ Old = \"51\u00d740=\"; New = \"88\u00d740=\" },\n    @{ Old = \"91\u00d711=\"; New = \"55\u00d769=\" },\n    @{ Old = \"18\u00d756=\"; New = \"21\u00d745=\" },\n    @{ Old = \"51\u00d742=\"; New = \"15\u00d785=\" },\n    @{ Old = \"32\u00d788=\"; New = \"56\u00d744=\" },\n    @{ Old = \"51\u00d784=\"; New = \"38\u00d790=\" },\n    @{ Old = \"95\u00d756=\"; New = \"97\u00d772=\" },\n    @{ Old = \"43\u00d757=\"; New = \"74\u00d735=\" },\n    @{ Old = \"77\u00d715=\"; New = \"20\u00d765=\" },\n    @{ Old = \"17\u00d720=\"; New = \"89\u00d737=\" },\n    @{ Old = \"91\u00d762=\"; New = \"86\u00d791=\" },\n    @{ Old = \"56\u00d723=\"; New = \"53\u00d744=\" },\n    @{ Old = \"38\u00d715=\"; New = \"64\u00d713=\" },\n    @{ Old = \"59\u00d749=\"; New = \"64\u00d712=\" },\n    @{ Old = \"46\u00d754=\"; New = \"29\u00d788=\" },\n    @{ Old = \"57\u00d785=\"; New = \"87\u00d738=\" },\n    @{ Old = \"80\u00d735=\"; New = \"34\u00d720=\" },\n    @{ Old = \"94\u00d781=\"; New = \"35\u00d764=\" },\n    @{ Old = \"77\u00d729=\"; New = \"41\u00d713=\" },\n    @{ Old = \"41\u00d719=\"; New = \"57\u00d711=\" },\n    @{ Old = \"43\u00d749=\"; New = \"11\u00d798=\" },\n    @{ Old = \"48\u00d712=\"; New = \"66\u00d725=\" },\n    @{ Old = \"40\u00d775=\"; New = \"62\u00d742=\" },\n    @{ Old = \"82\u00d744=\"; New = \"59\u00d743=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2) | Out-Null\n}\n"}
